# Add a new data row (row 21) to the "WorkSheet 1" sheet, mirroring the
# structure of the existing rows (dates in column A, numeric feature
# columns B:M, and a "Method" label in column N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 21

# Column A: date/time serial value, formatted like the other rows in the sheet
# (column A already carries the date/time style, so just set the value).
$ws.Cells.Item($row, 1).Value = 42625.883391203701

# Column B: feature value.
$ws.Cells.Item($row, 2).Value = 29

# Columns C through M: zeroed-out feature values.
for ($col = 3; $col -le 13; $col++) {
    $ws.Cells.Item($row, $col).Value = 0
}

# Column N: text label matching the other rows ("Random").
$ws.Cells.Item($row, 14).Value = "Random"
